$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings) ---
$ws.Range("A8").Value = "Volume 32   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/13/2025  Through  1/19/2025"

# --- Type-change cells: copy formats then copy values from a donor cell, or flip to numeric ---
$ws.Range("C14").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("F15").PasteSpecial(-4163)

$ws.Range("H16").Copy()
$ws.Range("N15").PasteSpecial(-4122)
$ws.Range("N15").Value = -100

$ws.Range("C14").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C17").PasteSpecial(-4163)

$ws.Range("C14").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C20").PasteSpecial(-4163)

$ws.Range("C14").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D20").PasteSpecial(-4163)

$ws.Range("E14").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E20").PasteSpecial(-4163)

$ws.Range("H16").Copy()
$ws.Range("L22").PasteSpecial(-4122)
$ws.Range("L22").Value = 50

$ws.Range("C14").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("F27").PasteSpecial(-4163)

$ws.Range("C14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("D28").PasteSpecial(-4163)

$ws.Range("E14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E28").PasteSpecial(-4163)

$excel.CutCopyMode = $false

# --- Simple numeric value updates (style unchanged) ---
$ws.Range("H15").Value = -100
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 3
$ws.Range("J16").Value = 6
$ws.Range("K16").Value = -50
$ws.Range("L16").Value = -76.923076923076
$ws.Range("M16").Value = -75
$ws.Range("N16").Value = -94
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -100
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = -11.111111111111
$ws.Range("J17").Value = 7
$ws.Range("K17").Value = -28.571428571428
$ws.Range("M17").Value = -16.666666666666
$ws.Range("N17").Value = -64.285714285714
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -20
$ws.Range("I18").Value = 8
$ws.Range("J18").Value = 14
$ws.Range("K18").Value = -42.857142857142
$ws.Range("L18").Value = -63.636363636363
$ws.Range("M18").Value = -55.555555555555
$ws.Range("N18").Value = -81.395348837209
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 22
$ws.Range("E19").Value = -45.454545454545
$ws.Range("F19").Value = 63
$ws.Range("G19").Value = 79
$ws.Range("H19").Value = -20.253164556962
$ws.Range("I19").Value = 50
$ws.Range("J19").Value = 53
$ws.Range("K19").Value = -5.66037735849
$ws.Range("L19").Value = -33.333333333333
$ws.Range("M19").Value = 13.636363636363
$ws.Range("N19").Value = -65.277777777777
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = -33.333333333333
$ws.Range("N20").Value = -97.674418604651
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = -53.125
$ws.Range("F21").Value = 95
$ws.Range("G21").Value = 121
$ws.Range("H21").Value = -21.487603305785
$ws.Range("I21").Value = 67
$ws.Range("J21").Value = 82
$ws.Range("K21").Value = -18.292682926829
$ws.Range("L21").Value = -43.697478991596
$ws.Range("M21").Value = -17.283950617283
$ws.Range("N21").Value = -77.28813559322
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 0
$ws.Range("I22").Value = 3
$ws.Range("J22").Value = 5
$ws.Range("K22").Value = -40
$ws.Range("M22").Value = -66.666666666666
$ws.Range("C24").Value = 33
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = 50
$ws.Range("F24").Value = 117
$ws.Range("G24").Value = 108
$ws.Range("H24").Value = 8.333333333333
$ws.Range("I24").Value = 81
$ws.Range("J24").Value = 64
$ws.Range("K24").Value = 26.5625
$ws.Range("L24").Value = -7.954545454545
$ws.Range("M24").Value = 22.727272727272
$ws.Range("C25").Value = 21
$ws.Range("D25").Value = 25
$ws.Range("E25").Value = -16
$ws.Range("F25").Value = 83
$ws.Range("G25").Value = 88
$ws.Range("H25").Value = -5.681818181818
$ws.Range("I25").Value = 53
$ws.Range("J25").Value = 51
$ws.Range("K25").Value = 3.92156862745
$ws.Range("L25").Value = -29.333333333333
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = -16.666666666666
$ws.Range("F26").Value = 17
$ws.Range("G26").Value = 32
$ws.Range("H26").Value = -46.875
$ws.Range("I26").Value = 14
$ws.Range("J26").Value = 21
$ws.Range("K26").Value = -33.333333333333
$ws.Range("L26").Value = -30
$ws.Range("M26").Value = 75
$ws.Range("H27").Value = -100
$ws.Range("C28").Value = 5
$ws.Range("F28").Value = 9
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 125
$ws.Range("I28").Value = 6
$ws.Range("K28").Value = 200
$ws.Range("L28").Value = 50
